$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Replace-ParaXml($para, [string]$innerXml) {
    # Replaces the full content (text + paragraph mark) of $para with the
    # supplied OOXML fragment (one or more <w:p> elements).
    $rng = $para.Range
    $rng.InsertXML($innerXml)
}

function Insert-ParaXmlBefore($para, [string]$innerXml) {
    $start = $para.Range.Start
    $insPoint = $d.Range($start, $start)
    $insPoint.InsertXML($innerXml)
}

# ---------------------------------------------------------------------
# 1. New "Submit paper" paragraph + blank paragraph before "Intro:"
# ---------------------------------------------------------------------
$introPara = $d.Paragraphs.Item(1)
$xml1 = "<w:p $wns><w:r><w:t>Submit paper</w:t></w:r></w:p><w:p/>"
Insert-ParaXmlBefore $introPara $xml1

Write-Host "Step1 count:" $d.Paragraphs.Count
Write-Host "Step1 p1:" $d.Paragraphs.Item(1).Range.Text
Write-Host "Step1 p2:" $d.Paragraphs.Item(2).Range.Text
Write-Host "Step1 p3:" $d.Paragraphs.Item(3).Range.Text

# ---------------------------------------------------------------------
# 2. "Tell about LIP/ICP first" -> "Check what kind of referencing
#    should be used"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Tell about LIP/ICP first", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Check what kind of referencing should be used", 2) | Out-Null

Write-Host "Step2 p4:" $d.Paragraphs.Item(4).Range.Text

# ---------------------------------------------------------------------
# 3. Remove "Fuse flywheel model and pendulums? ..." and "Write related
#    works" paragraphs entirely, and remove the text (not the paragraph
#    mark) of "Theoretic Limits:" so its mark -- and the bookmark that
#    lives on it -- merges into the "Background:" paragraph.
# ---------------------------------------------------------------------
$fusePara = $d.Paragraphs.Item(6)
$writeRelatedPara = $d.Paragraphs.Item(7)
$rngFW = $d.Range($fusePara.Range.Start, $writeRelatedPara.Range.End)
$rngFW.Delete() | Out-Null

$theoLimitsPara = $d.Paragraphs.Item(6)   # now "Theoretic Limits:"
Write-Host "Step3 theoLimits text:" $theoLimitsPara.Range.Text
$rngText = $d.Range($theoLimitsPara.Range.Start, $theoLimitsPara.Range.End - 1)
$rngText.Delete() | Out-Null

$backgroundPara = $d.Paragraphs.Item(5)
$rngMark = $d.Range($backgroundPara.Range.End - 1, $backgroundPara.Range.End)
$rngMark.Delete() | Out-Null

Write-Host "Step3 count:" $d.Paragraphs.Count
for ($i=1; $i -le 9; $i++) {
  Write-Host "Step3" $i ":" $d.Paragraphs.Item($i).Range.Text
}

# ---------------------------------------------------------------------
# 4. "Make unilateral and height constrained concise and write proof
#    separate." -> "Change flywheel model in: 'dynamics'"
# ---------------------------------------------------------------------
$d.Paragraphs.Item(6).Range.Find.Execute(
    "Make unilateral and height constrained concise and write proof separate.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Change flywheel model in: ‘dynamics’", 2) | Out-Null

Write-Host "Step4 p6:" $d.Paragraphs.Item(6).Range.Text

# ---------------------------------------------------------------------
# 5. "Add force constrained limits" -> "Incorporate that CMP is not
#    restricted to lie inside the polygon"
# ---------------------------------------------------------------------
$d.Paragraphs.Item(7).Range.Find.Execute(
    "Add force constrained limits", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Incorporate that CMP is not restricted to lie inside the polygon", 2) | Out-Null

Write-Host "Step5 p7:" $d.Paragraphs.Item(7).Range.Text

# ---------------------------------------------------------------------
# 6. Replace "Extend capturability comparison -> fix comparison plot"
#    with two new paragraphs: "Theoretic Limits:" and the list item
#    "Write about velocity plot".
# ---------------------------------------------------------------------
$extendPara = $d.Paragraphs.Item(8)
Write-Host "Step6 before:" $extendPara.Range.Text
$xml6 = "<w:p $wns><w:r><w:t>Theoretic Limits:</w:t></w:r></w:p>" +
        "<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr>" +
        "<w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
        "<w:r><w:t>Write about velocity plot</w:t></w:r></w:p>"
$extendPara.Range.InsertXML($xml6)

Write-Host "Step6 count:" $d.Paragraphs.Count
for ($i=1; $i -le 12; $i++) {
  Write-Host "Step6" $i ":" $d.Paragraphs.Item($i).Range.Text
}

# ---------------------------------------------------------------------
# 7. Remove "Rewrite a bit and write discussion why not." paragraph.
# ---------------------------------------------------------------------
$rewritePara = $d.Paragraphs.Item(11)
Write-Host "Step7 remove:" $rewritePara.Range.Text
$rewritePara.Range.Delete() | Out-Null

Write-Host "Step7 count:" $d.Paragraphs.Count
for ($i=9; $i -le 13; $i++) {
  Write-Host "Step7" $i ":" $d.Paragraphs.Item($i).Range.Text
}

# ---------------------------------------------------------------------
# 8. Remove "Run 360 push" paragraph.
# ---------------------------------------------------------------------
$run360Para = $d.Paragraphs.Item(12)
Write-Host "Step8 remove:" $run360Para.Range.Text
$run360Para.Range.Delete() | Out-Null

Write-Host "Step8 count:" $d.Paragraphs.Count
for ($i=10; $i -le 14; $i++) {
  Write-Host "Step8" $i ":" $d.Paragraphs.Item($i).Range.Text
}

# ---------------------------------------------------------------------
# 9. "Try some angular momentum tests." -> "Try some angular momentum
#    tests 360 push"
# ---------------------------------------------------------------------
$d.Paragraphs.Item(12).Range.Find.Execute(
    "Try some angular momentum tests.", $true, $false, $false, $false,
    $false, $true, 1, $false, "Try some angular momentum tests 360 push", 2) | Out-Null

Write-Host "Step9 p12:" $d.Paragraphs.Item(12).Range.Text

# ---------------------------------------------------------------------
# 10 & 11. Remove "Hardware: pd gains foot angular " and "Ask feedback
#          paper" paragraphs (both entirely).
# ---------------------------------------------------------------------
$hardwarePara = $d.Paragraphs.Item(13)
$askFeedbackPara = $d.Paragraphs.Item(14)
Write-Host "Step10 remove:" $hardwarePara.Range.Text "|" $askFeedbackPara.Range.Text
$rngHW = $d.Range($hardwarePara.Range.Start, $askFeedbackPara.Range.End)
$rngHW.Delete() | Out-Null

Write-Host "Step10 count:" $d.Paragraphs.Count
for ($i=11; $i -le 16; $i++) {
  Write-Host "Step10" $i ":" $d.Paragraphs.Item($i).Range.Text
}

# ---------------------------------------------------------------------
# 12. "Atlas hardware / ball" gains three more runs:
#       bold " -> check "; "pd" (spell-check wrapped); " gains foot
#       angular"
# ---------------------------------------------------------------------
$atlasPara = $d.Paragraphs.Item(13)
Write-Host "Step12 before:" $atlasPara.Range.Text
$xml12 = "<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr>" +
         "<w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
         "<w:r><w:rPr><w:b/></w:rPr><w:t>Atlas hardware / ball</w:t></w:r>" +
         "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'> -> check </w:t></w:r>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>pd</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "<w:r><w:t xml:space='preserve'> gains foot angular</w:t></w:r></w:p>"
$atlasPara.Range.InsertXML($xml12)

Write-Host "Step12 after:" $d.Paragraphs.Item(13).Range.Text
Write-Host "Step12 count:" $d.Paragraphs.Count

# ---------------------------------------------------------------------
# 13. "Run 360 push incremental 0.0 0.1 … 0.7 swing phase" gains a
#     new run: " -> fix after 0.5"
# ---------------------------------------------------------------------
$run360IncPara = $d.Paragraphs.Item(17)
Write-Host "Step13 before:" $run360IncPara.Range.Text
$xml13 = "<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr>" +
         "<w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
         "<w:r><w:t>Run 360 push incremental 0.0 0.1 … 0.7 swing phase</w:t></w:r>" +
         "<w:r><w:t xml:space='preserve'> -> fix after 0.5</w:t></w:r></w:p>"
$run360IncPara.Range.InsertXML($xml13)

Write-Host "Step13 after:" $d.Paragraphs.Item(17).Range.Text
Write-Host "Step13 count:" $d.Paragraphs.Count

# ---------------------------------------------------------------------
# 14. The trailing empty list-item paragraph loses its list formatting
#     and becomes "<tab>-"; then a batch of new paragraphs (Bibliography
#     section + closing remark) is appended after it.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
Write-Host "Step14 before text:" $lastPara.Range.Text
$xml14 = "<w:p $wns><w:r><w:tab/><w:t>-</w:t></w:r></w:p>" +
         "<w:p><w:r><w:t>Bibliography</w:t></w:r></w:p>" +
         "<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr>" +
         "<w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
         "<w:r><w:t>Entries?</w:t></w:r></w:p>" +
         "<w:p><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr>" +
         "<w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
         "<w:r><w:t>Citing and citing of names</w:t></w:r>" +
         "<w:r><w:t xml:space='preserve'> in thesis?</w:t></w:r></w:p>" +
         "<w:p/>" +
         "<w:p><w:r><w:t>Focus on melting everything together</w:t></w:r></w:p>" +
         "<w:p/>"
$lastPara.Range.InsertXML($xml14)

Write-Host "Step14 count:" $d.Paragraphs.Count
for ($i=19; $i -le $d.Paragraphs.Count; $i++) {
  Write-Host "Step14" $i ":" $d.Paragraphs.Item($i).Range.Text
}
